# ToDo workbook edit
# - Remove Sheet2 and Sheet3, rename Sheet1 to "ToDo"
# - Add "Wer?" (C), "Erledigt?" (D) and "Problem" (E) columns to the ToDo list
# - Add four new rows of data (12-15)
# - Center align the new Wer?/Erledigt? data cells
# - Re-apply an AutoFilter over A1:D15 that hides rows whose "Erledigt?" column is not blank

$excel.DisplayAlerts = $false

$wb = $excel.ActiveWorkbook

# --- remove the unused sheets, keep only the ToDo list ---------------------
$wb.Worksheets("Sheet2").Delete()
$wb.Worksheets("Sheet3").Delete()

$ws = $wb.Worksheets("Sheet1")
$ws.Name = "ToDo"

# --- drop the old autofilter so the new one can take the full range --------
$ws.AutoFilterMode = $false

# --- header row --------------------------------------------------------------
$ws.Range("C1").Value = "Wer?"
$ws.Range("D1").Value = "Erledigt?"
$ws.Range("E1").Value = "Problem"

# --- existing rows: who worked on it / done flag ----------------------------
$ws.Range("C3").Value = "HS"

$ws.Range("C4").Value = "PH"
$ws.Range("D4").Value = "X"

$ws.Range("C5").Value = "HS"
$ws.Range("C6").Value = "HS"

$ws.Range("C8").Value = "PH"
$ws.Range("D8").Value = "X"

$ws.Range("C9").Value = "PH"
$ws.Range("C10").Value = "HS"
$ws.Range("C11").Value = "PH"

# extra "Problem" notes for row 11
$ws.Range("E11").Value = "<h2> kleinere Schriftgröße klappt nicht"
$ws.Range("F11").Value = "F12 im Chrome"

# --- new rows 12-15 ----------------------------------------------------------
$ws.Range("A12").Value = "Kategorien"
$ws.Range("B12").Value = "Bei Neu: Überkategorie kann nicht ausgewählt werden"
$ws.Range("C12").Value = "HS"

$ws.Range("A13").Value = "Allgemein"
$ws.Range("B13").Value = "Alle Warnings und Fehler nicht auf der Seite anzeigen"

$ws.Range("A14").Value = "Allgemein"
$ws.Range("B14").Value = "Datenbank online bringen"

$ws.Range("A15").Value = "Allgemein"
$ws.Range("B15").Value = "Footer entfernen"
$ws.Range("C15").Value = "PH"
$ws.Range("D15").Value = "X"

# --- center align every "Wer?"/"Erledigt?" data cell ------------------------
# rows 2-12 and 15 have both a "Wer?" and "Erledigt?" cell, rows 13-14 only
# carry the (empty) "Erledigt?" cell
$both1 = $ws.Range("C2:D12")
$both1.HorizontalAlignment = -4108
$both1.VerticalAlignment = -4108

$both2 = $ws.Range("C15:D15")
$both2.HorizontalAlignment = -4108
$both2.VerticalAlignment = -4108

$doneOnly = $ws.Range("D13:D14")
$doneOnly.HorizontalAlignment = -4108
$doneOnly.VerticalAlignment = -4108

# --- selection / active cell mirrors the author's saved view ----------------
$ws.Range("A16").Select()

# --- autofilter: show only rows where "Erledigt?" (col D) is blank ----------
$ws.Range("A1:D15").AutoFilter(4, "", 7)
